$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("F2").Value = 1.7
$ws.Range("G2").Value = 1.71
$ws.Range("I2").Value = 5.6
$ws.Range("J2").Value = 4.2
$ws.Range("K2").Value = 4.3
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 1.43
$ws.Range("S2").Value = 3.15
$ws.Range("W2").Value = 2.4
$ws.Range("Y2").Value = 20
$ws.Range("Z2").Value = 44
$ws.Range("AA2").Value = 140
$ws.Range("AD2").Value = 21
$ws.Range("AE2").Value = 70
$ws.Range("AH2").Value = 20
$ws.Range("AN2").Value = 9.6
$ws.Range("AO2").Value = 75

# Row 4 updates
$ws.Range("F4").Value = 8.4
$ws.Range("G4").Value = 9.4
$ws.Range("H4").Value = 1.41
$ws.Range("I4").Value = 1.43
$ws.Range("N4").Value = 4
$ws.Range("Q4").Value = 1.87
$ws.Range("S4").Value = 3.2
$ws.Range("T4").Value = 2.22
$ws.Range("V4").Value = 3.3
$ws.Range("AL4").Value = 150
